# SectorGroup.xlsx column reorder:
#   codeforiati:category-name <-> codeforiati:group-name   (column D <-> column E)
#   codeforiati:group-code    <-> codeforiati:category-code (column F <-> column G)
# The category/group names now lead with the broader "group" value before the
# more specific "category" value, matching the upstream codeforIATI/codelists
# export ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$rangeD = $ws.Range("D1:D$lastRow")
$rangeE = $ws.Range("E1:E$lastRow")
$rangeF = $ws.Range("F1:F$lastRow")
$rangeG = $ws.Range("G1:G$lastRow")

# Scratch columns, well outside the used A:G range, used purely to hold a
# value mid-swap (Copy preserves the source's text/number type so the
# shared-string "t=s" typing of the code columns survives the round-trip).
$tempD = $ws.Range("Z1:Z$lastRow")
$tempF = $ws.Range("AA1:AA$lastRow")

$rangeD.Copy($tempD)
$rangeE.Copy($rangeD)
$tempD.Copy($rangeE)

$rangeF.Copy($tempF)
$rangeG.Copy($rangeF)
$tempF.Copy($rangeG)

$ws.Range("Z1:AA$lastRow").ClearContents()
